# Improve responsive design implementation (resolves issue #31)
#
# Appends one new data row (row 45) to each of the four worksheets, right
# after the existing last row (row 44), carrying the same column layout:
#   A: timestamp (date/time, formatted like the existing rows)
#   B-E: hex-dump text fields (kept as literal text)
#   F-I: numeric fields
#
# The four sheets are processed in tab order (1..4), matching
# ROW35-FE-LIFTER, ROW35-MID-LIFTER, ROW02-FE-LIFTER, ROW02-MID-LIFTER.

$wb = $excel.ActiveWorkbook

$newRows = @(
    @{
        Sheet = 1
        A = [double]"45746.83676112269"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x76"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 374
        I = 13
    },
    @{
        Sheet = 2
        A = [double]"45746.68584070602"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x76"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 374
        I = 14
    },
    @{
        Sheet = 3
        A = [double]"45746.82797267361"
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x76"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 374
        I = 3
    },
    @{
        Sheet = 4
        A = [double]"45746.88490554398"
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x76"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 374
        I = 3
    }
)

foreach ($row in $newRows) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    # New row goes right below the current last row (44 -> 45).
    $targetRow = 45

    $ws.Cells.Item($targetRow, 1).Value = $row.A
    $ws.Cells.Item($targetRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($targetRow, 2).Value = $row.B
    $ws.Cells.Item($targetRow, 3).Value = $row.C
    $ws.Cells.Item($targetRow, 4).Value = $row.D
    $ws.Cells.Item($targetRow, 5).Value = $row.E
    $ws.Cells.Item($targetRow, 6).Value = $row.F
    $ws.Cells.Item($targetRow, 7).Value = $row.G
    $ws.Cells.Item($targetRow, 8).Value = $row.H
    $ws.Cells.Item($targetRow, 9).Value = $row.I
}
